# Update "想去人数" (number of people interested) figures that changed
# between the two data refreshes captured in the commit.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 265
$wsExhibit.Range("F4").Value = 915

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 265
$wsAll.Range("F5").Value = 915
